$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '29.036.91'
Set-TextValue 'E2' '  -0.66%  '
Set-TextValue 'D3' '1.830.56'
Set-TextValue 'E3' '  -0.71%  '
Set-TextValue 'D4' '0.9995'
Set-TextValue 'E4' '  +0.05%  '
Set-TextValue 'D5' '241.73'
Set-TextValue 'E5' '  +0.43%  '
Set-TextValue 'D6' '0.6539'
Set-TextValue 'E6' '  -2.95%  '
Set-TextValue 'D7' '1.001'
Set-TextValue 'D8' '44.51'
Set-TextValue 'E8' '  +5.84%  '
Set-TextValue 'D9' '0.2934'
Set-TextValue 'E9' '  -0.66%  '
Set-TextValue 'D10' '0.07328'
Set-TextValue 'E10' '  -1.31%  '
Set-TextValue 'D11' '22.91'
Set-TextValue 'E11' '  +0.15%  '
Set-TextValue 'D12' '0.07662'
Set-TextValue 'E12' '  -0.70%  '
Set-TextValue 'D13' '1.837.51'
Set-TextValue 'E13' '  -0.09%  '
Set-TextValue 'D14' '4.977'
Set-TextValue 'E14' '  -0.66%  '
Set-TextValue 'D15' '0.6664'
Set-TextValue 'E15' '  -0.72%  '
Set-TextValue 'D16' '81.68'
Set-TextValue 'E16' '  -5.14%  '
Set-TextValue 'D17' '6.101'
Set-TextValue 'E17' '  -0.48%  '
Set-TextValue 'D18' '0.000008694'
Set-TextValue 'E18' '  +4.59%  '
Set-TextValue 'D19' '29.041.13'
Set-TextValue 'E19' '  -0.50%  '
Set-TextValue 'D20' '2.087.88'
Set-TextValue 'E20' '  +0.42%  '
Set-TextValue 'D21' '12.43'
Set-TextValue 'E21' '  -0.69%  '
Set-TextValue 'D22' '223.44'
Set-TextValue 'E22' '  -2.24%  '
Set-TextValue 'E23' '  +0.04%  '
Set-TextValue 'D24' '7.111'
Set-TextValue 'E24' '  -1.07%  '
Set-TextValue 'D26' '157.51'
Set-TextValue 'E26' '  -1.96%  '
Set-TextValue 'D27' '8.482'
Set-TextValue 'E27' '  -2.43%  '
Set-TextValue 'D28' '0.1377'
Set-TextValue 'E28' '  -1.90%  '
Set-TextValue 'D29' '17.90'
Set-TextValue 'E29' '  -0.56%  '
Set-TextValue 'D30' '1.502'
Set-TextValue 'E30' '  -0.41%  '
Set-TextValue 'D31' '4.101'
Set-TextValue 'E31' '  -1.81%  '
Set-TextValue 'D32' '4.010'
Set-TextValue 'E32' '  -1.42%  '
Set-TextValue 'D33' '1.199'
Set-TextValue 'E33' '  +0.74%  '
Set-TextValue 'D34' '0.05345'
Set-TextValue 'E34' '  +0.69%  '
Set-TextValue 'D35' '1.837'
Set-TextValue 'E35' '  -2.01%  '
Set-TextValue 'D36' '0.7418'
Set-TextValue 'E36' '  -2.40%  '
Set-TextValue 'D37' '1.156'
Set-TextValue 'E37' '  +1.77%  '
Set-TextValue 'D38' '2.650'
Set-TextValue 'E38' '  -0.98%  '
Set-TextValue 'D39' '1.297.16'
Set-TextValue 'E39' '  -2.38%  '
Set-TextValue 'D40' '0.01787'
Set-TextValue 'E40' '  -0.97%  '
Set-TextValue 'D41' '2.748'
Set-TextValue 'E41' '  +0.65%  '
Set-TextValue 'D42' '6.320'
Set-TextValue 'E42' '  +6.24%  '
Set-TextValue 'D43' '0.8992'
Set-TextValue 'E43' '  -2.11%  '
Set-TextValue 'D44' '0.9999'
Set-TextValue 'D45' '103.29'
Set-TextValue 'E45' '  -0.13%  '
Set-TextValue 'B46' 'RocketPoolETH'
Set-TextValue 'C46' 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextValue 'D46' '1.986.80'
Set-TextValue 'E46' '  +0.46%  '
Set-TextValue 'B47' 'XinFinNetwork'
Set-TextValue 'C47' 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
Set-TextValue 'D47' '0.07857'
Set-TextValue 'E47' '  -2.24%  '
Set-TextValue 'D48' '64.33'
Set-TextValue 'E48' '  +0.77%  '
Set-TextValue 'D49' '0.5133'
Set-TextValue 'E49' '  -0.56%  '
Set-TextValue 'D50' '0.00000000120'
Set-TextValue 'E50' '  -1.82%  '
Set-TextValue 'D51' '1.736'
Set-TextValue 'E51' '  -2.14%  '
